# Auto-generated Excel COM-interop script to apply Ifrit_Profits.xlsx diff
# Updates currentAveragePrice/LevePrice/LeveProfit columns (H:N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1993.4706
$ws.Range("I32").Value = 2516.875
$ws.Range("J32").Value = 1528.2222
$ws.Range("K32").Value = 2516.875
$ws.Range("L32").Value = 1528.2222
$ws.Range("M32").Value = -2190.875
$ws.Range("N32").Value = -2180.2222
$ws.Range("H51").Value = 4997
$ws.Range("J51").Value = 4997
$ws.Range("L51").Value = 4997
$ws.Range("N51").Value = -5965
$ws.Range("H125").Value = 950.2857
$ws.Range("I125").Value = 288
$ws.Range("J125").Value = 1833.3334
$ws.Range("K125").Value = 2592
$ws.Range("L125").Value = 16500.0006
$ws.Range("M125").Value = -132
$ws.Range("N125").Value = -21420.0006
$ws.Range("H132").Value = 532565.5
$ws.Range("I132").Value = 532565.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1597696.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1595166.5
$ws.Range("N132").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11935.162
$ws.Range("I32").Value = 5692.32
$ws.Range("J32").Value = 24941.084
$ws.Range("K32").Value = 5692.32
$ws.Range("L32").Value = 24941.084
$ws.Range("M32").Value = -5405.32
$ws.Range("N32").Value = -25515.084

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2668.6316
$ws.Range("I107").Value = 2161.6667
$ws.Range("J107").Value = 3124.9
$ws.Range("K107").Value = 2161.6667
$ws.Range("L107").Value = 3124.9
$ws.Range("M107").Value = -241.6667000000002
$ws.Range("N107").Value = -6964.9
$ws.Range("H134").Value = 31909.111
$ws.Range("I134").Value = 34537.816
$ws.Range("J134").Value = 2993.3333
$ws.Range("K134").Value = 103613.448
$ws.Range("L134").Value = 8979.999899999999
$ws.Range("M134").Value = -101078.448
$ws.Range("N134").Value = -14049.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 919.0769
$ws.Range("I105").Value = 768.1818
$ws.Range("J105").Value = 1749
$ws.Range("K105").Value = 768.1818
$ws.Range("L105").Value = 1749
$ws.Range("M105").Value = 978.8182
$ws.Range("N105").Value = -5243
$ws.Range("H107").Value = 2321.6667
$ws.Range("I107").Value = 2279.25
$ws.Range("J107").Value = 2406.5
$ws.Range("K107").Value = 2279.25
$ws.Range("L107").Value = 2406.5
$ws.Range("M107").Value = -359.25
$ws.Range("N107").Value = -6246.5
$ws.Range("H122").Value = 17857892
$ws.Range("I122").Value = 17857892
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 53573676
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -53571226
$ws.Range("N122").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1141.4286
$ws.Range("I3").Value = 831.6667
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 2495.0001
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -2383.0001
$ws.Range("N3").Value = -9224
$ws.Range("H60").Value = 27779336
$ws.Range("I60").Value = 47619504
$ws.Range("J60").Value = 3100
$ws.Range("K60").Value = 142858512
$ws.Range("L60").Value = 9300
$ws.Range("M60").Value = -142858261
$ws.Range("N60").Value = -9802
$ws.Range("H129").Value = 1400
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1400
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 4200
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -14200
$ws.Range("H133").Value = 6955
$ws.Range("I133").Value = 5200
$ws.Range("J133").Value = 7058.2354
$ws.Range("K133").Value = 15600
$ws.Range("L133").Value = 21174.7062
$ws.Range("M133").Value = -10540
$ws.Range("N133").Value = -31294.7062
$ws.Range("H134").Value = 38465252
$ws.Range("I134").Value = 38465252
$ws.Range("K134").Value = 115395756
$ws.Range("M134").Value = -115390686
$ws.Range("H136").Value = 2343
$ws.Range("I136").Value = 2343
$ws.Range("K136").Value = 7029
$ws.Range("M136").Value = -1929
$ws.Range("H138").Value = 3162.5
$ws.Range("I138").Value = 3162.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9487.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -4347.5
$ws.Range("N138").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1834
$ws.Range("I7").Value = 1834
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1834
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1722
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 671.55554
$ws.Range("I22").Value = 737.7143
$ws.Range("J22").Value = 440
$ws.Range("K22").Value = 737.7143
$ws.Range("L22").Value = 440
$ws.Range("M22").Value = -442.7143
$ws.Range("N22").Value = -1030
$ws.Range("H27").Value = 671.55554
$ws.Range("I27").Value = 737.7143
$ws.Range("J27").Value = 440
$ws.Range("K27").Value = 737.7143
$ws.Range("L27").Value = 440
$ws.Range("M27").Value = -630.7143
$ws.Range("N27").Value = -654
$ws.Range("H61").Value = 3383.3333
$ws.Range("I61").Value = 2060
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 2060
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -1858
$ws.Range("N61").Value = -10404
$ws.Range("H113").Value = 3383.3333
$ws.Range("I113").Value = 2060
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 2060
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 110
$ws.Range("N113").Value = -14340
$ws.Range("H126").Value = 1834
$ws.Range("I126").Value = 1834
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5502
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3032
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4743.7905
$ws.Range("I132").Value = 5082.029
$ws.Range("J132").Value = 3264
$ws.Range("K132").Value = 15246.087
$ws.Range("L132").Value = 9792
$ws.Range("M132").Value = -12716.087
$ws.Range("N132").Value = -14852
$ws.Range("H136").Value = 1948.8
$ws.Range("I136").Value = 1320
$ws.Range("J136").Value = 4014.8572
$ws.Range("K136").Value = 3960
$ws.Range("L136").Value = 12044.5716
$ws.Range("M136").Value = -1410
$ws.Range("N136").Value = -17144.5716

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1054.8948
$ws.Range("I126").Value = 943.7059
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 2831.1177
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -361.1177000000002
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2495.3513
$ws.Range("I132").Value = 2648
$ws.Range("K132").Value = 7944
$ws.Range("M132").Value = -5414
